$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.228476821192053
$ws.Range("C2").Value = 0.5
$ws.Range("J2").Value = 0.01655629139072848
$ws.Range("P2").Value = 0.1821192052980132
$ws.Range("S2").Value = 0.0728476821192053

# Row 3
$ws.Range("C3").Value = 0.03773584905660377
$ws.Range("J3").Value = 0.03773584905660377
$ws.Range("P3").Value = 0.7924528301886793
$ws.Range("S3").Value = 0.1320754716981132

# Row 4
$ws.Range("J4").Value = 0.1388888888888889
$ws.Range("P4").Value = 0.6388888888888888
$ws.Range("S4").Value = 0.2222222222222222

# Row 6
$ws.Range("B6").Value = 0.0599250936329588
$ws.Range("D6").Value = 0.00749063670411985
$ws.Range("F6").Value = 0.0898876404494382
$ws.Range("J6").Value = 0.2359550561797753
$ws.Range("O6").Value = 0.02621722846441948
$ws.Range("Q6").Value = 0.1385767790262172
$ws.Range("R6").Value = 0.09737827715355805
$ws.Range("S6").Value = 0.3445692883895131

# Row 7
$ws.Range("B7").Value = 0.1013215859030837
$ws.Range("D7").Value = 0.01762114537444934
$ws.Range("E7").Value = 0.004405286343612335
$ws.Range("F7").Value = 0.1145374449339207
$ws.Range("J7").Value = 0.1497797356828194
$ws.Range("O7").Value = 0.02202643171806168
$ws.Range("Q7").Value = 0.1718061674008811
$ws.Range("R7").Value = 0.08370044052863436
$ws.Range("S7").Value = 0.3348017621145374

# Row 8
$ws.Range("B8").Value = 0.09200968523002422
$ws.Range("D8").Value = 0.01210653753026634
$ws.Range("F8").Value = 0.08716707021791767
$ws.Range("J8").Value = 0.1089588377723971
$ws.Range("O8").Value = 0.02179176755447942
$ws.Range("Q8").Value = 0.1646489104116223
$ws.Range("R8").Value = 0.08958837772397095
$ws.Range("S8").Value = 0.423728813559322

# Row 9
$ws.Range("B9").Value = 0.09036144578313253
$ws.Range("D9").Value = 0.01204819277108434
$ws.Range("F9").Value = 0.05421686746987952
$ws.Range("J9").Value = 0.1265060240963855
$ws.Range("O9").Value = 0.03614457831325301
$ws.Range("Q9").Value = 0.1566265060240964
$ws.Range("R9").Value = 0.1325301204819277
$ws.Range("S9").Value = 0.3915662650602409

# Row 10
$ws.Range("B10").Value = 0.1085450346420323
$ws.Range("D10").Value = 0.01770592763664357
$ws.Range("F10").Value = 0.07082371054657428
$ws.Range("J10").Value = 0.1293302540415704
$ws.Range("O10").Value = 0.0207852193995381
$ws.Range("Q10").Value = 0.2001539645881447
$ws.Range("R10").Value = 0.09083910700538876
$ws.Range("S10").Value = 0.3618167821401078

# Row 11
$ws.Range("G11").Value = 0.1311475409836066
$ws.Range("J11").Value = 0.09016393442622951
$ws.Range("K11").Value = 0.1775956284153005
$ws.Range("L11").Value = 0.5819672131147541
$ws.Range("S11").Value = 0.01912568306010929

# Row 12
$ws.Range("G12").Value = 0.6866359447004609
$ws.Range("J12").Value = 0.2488479262672811
$ws.Range("K12").Value = 0.009216589861751152
$ws.Range("L12").Value = 0.03225806451612903
$ws.Range("S12").Value = 0.02304147465437788

# Row 13
$ws.Range("G13").Value = 0.6545454545454545
$ws.Range("J13").Value = 0.3090909090909091
$ws.Range("S13").Value = 0.03636363636363636

# Row 14
$ws.Range("G14").Value = 0.6
$ws.Range("J14").Value = 0.4

# Row 15
$ws.Range("F15").Value = 0.03375527426160337
$ws.Range("H15").Value = 0.1687763713080169
$ws.Range("I15").Value = 0.02531645569620253
$ws.Range("J15").Value = 0.29957805907173
$ws.Range("K15").Value = 0.1265822784810127
$ws.Range("M15").Value = 0.02531645569620253
$ws.Range("O15").Value = 0.06751054852320675
$ws.Range("S15").Value = 0.2531645569620253

# Row 16
$ws.Range("F16").Value = 0.03980099502487562
$ws.Range("H16").Value = 0.1044776119402985
$ws.Range("I16").Value = 0.05970149253731343
$ws.Range("J16").Value = 0.4577114427860697
$ws.Range("K16").Value = 0.1243781094527363
$ws.Range("M16").Value = 0.01492537313432836
$ws.Range("O16").Value = 0.05472636815920398
$ws.Range("S16").Value = 0.1442786069651741

# Row 17
$ws.Range("F17").Value = 0.01605504587155963
$ws.Range("H17").Value = 0.1720183486238532
$ws.Range("I17").Value = 0.08944954128440367
$ws.Range("J17").Value = 0.4036697247706422
$ws.Range("K17").Value = 0.1238532110091743
$ws.Range("M17").Value = 0.02752293577981652
$ws.Range("N17").Value = 0.002293577981651376
$ws.Range("O17").Value = 0.06192660550458716
$ws.Range("S17").Value = 0.1032110091743119

# Row 18
$ws.Range("F18").Value = 0.009049773755656109
$ws.Range("H18").Value = 0.1809954751131222
$ws.Range("I18").Value = 0.08597285067873303
$ws.Range("J18").Value = 0.3936651583710407
$ws.Range("K18").Value = 0.1085972850678733
$ws.Range("M18").Value = 0.03167420814479638
$ws.Range("O18").Value = 0.07692307692307693
$ws.Range("S18").Value = 0.1131221719457014

# Row 19
$ws.Range("F19").Value = 0.01866883116883117
$ws.Range("H19").Value = 0.1964285714285714
$ws.Range("I19").Value = 0.07142857142857142
$ws.Range("J19").Value = 0.3603896103896104
$ws.Range("K19").Value = 0.1290584415584416
$ws.Range("M19").Value = 0.02353896103896104
$ws.Range("N19").Value = 0.003246753246753247
$ws.Range("O19").Value = 0.07061688311688312
$ws.Range("S19").Value = 0.1266233766233766
